# WIP cleanup pass: trim the stale trailing rows left on "fundamental_data",
# rename "projected_ei_in_Wh" to "projected_ei", and leave the
# "projected_ei" tab as the active/selected one (it was "historic_data"
# before).

$wb = $excel.ActiveWorkbook

# --- 1. fundamental_data: drop the leftover empty/placeholder rows 32-56 ---
# (only column P had stray styled-but-empty cells there; the real data ends
# at row 31). Deleting the rows also shrinks the sheet dimension and the
# Table1 range/autofilter automatically.
$fund = $wb.Worksheets.Item("fundamental_data")
$fund.Rows.Item(32).Resize(25).Delete()

# --- 2. rename the sheet ---
$proj = $wb.Worksheets.Item("projected_ei_in_Wh")
$proj.Name = "projected_ei"

# --- 3. move the active tab from historic_data to projected_ei ---
$proj.Activate()
